# Updated cryptos list with GitHub Actions
#
# Column D holds price strings that look numeric (e.g. "1.00", "10.94").
# In the source workbook these are plain text cells (t="inlineStr"), so we
# force column D writes through a Text number format and then restore the
# "Normal" style so the cell keeps its original (default) styling while the
# stored value remains an exact text string instead of being auto-converted
# to a number by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.391.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.904.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.901.27"
$ws.Range("D7").Style = "Normal"

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("E10").Value = "  -2.25%  "

$ws.Range("E11").Value = "  -0.77%  "

$ws.Range("E12").Value = "  -0.96%  "

$ws.Range("E13").Value = "  +2.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.561.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.910.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.498.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.27%  "

$ws.Range("E20").Value = "  -0.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "469.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.743"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000161"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.056.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.34%  "

$ws.Range("E33").Value = "  -2.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.875.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.12%  "

$ws.Range("E37").Value = "  -1.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.75"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +14.77%  "

$ws.Range("E39").Value = "  -0.58%  "

$ws.Range("E40").Value = "  +0.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("E43").Value = "  -1.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000307"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "426.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.42%  "

$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "28.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.97%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "47.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.60%  "
